$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ridership")

$ws.Range("C2").Value = 185
$ws.Range("D2").Value = 97.46

$ws.Range("C3").Value = 198
$ws.Range("D3").Value = 99.75

$ws.Range("C4").Value = 191
$ws.Range("D4").Value = 105.84

$ws.Range("C5").Value = 206
$ws.Range("D5").Value = 104.24

$ws.Range("C6").Value = 189
$ws.Range("D6").Value = 101.79

$ws.Range("C7").Value = 303
$ws.Range("D7").Value = 47.9

$ws.Range("C8").Value = 118
$ws.Range("D8").Value = 36.66
